# Generate Report for handback
# Update the "Latest Handoff Datetime" and "Latest Handback DateTime" values
# for the first source file row (100a6e02-...) on both the zh-cn and de-de
# language sheets, simulating a freshly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-08 10:19:59"
$wsZhCn.Range("G2").Value = "2016-01-08 10:20:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-08 10:20:10"
$wsDeDe.Range("G2").Value = "2016-01-08 10:21:00"
